$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 87) with the latest Argent (silver) price data.
# Values are entered as text (quote-prefixed) to match the existing sheet's
# convention of storing all price figures as text rather than numbers, then
# the style is reset to "Normal" so no extra text-format styling is applied.
$row = 87
$ws.Cells.Item($row, 1).Formula  = "'2025-05-27"
$ws.Cells.Item($row, 2).Formula  = "'35.5"
$ws.Cells.Item($row, 3).Formula  = "'35.11"
$ws.Cells.Item($row, 4).Formula  = "'0.94"
$ws.Cells.Item($row, 5).Formula  = "'0.258"
$ws.Cells.Item($row, 6).Formula  = "'0.09"
$ws.Cells.Item($row, 7).Formula  = "'5,386"
$ws.Cells.Item($row, 8).Formula  = "'8,063"
$ws.Cells.Item($row, 9).Formula  = "'8,113"
$ws.Cells.Item($row, 10).Formula = "'7.1925"

$ws.Range("A$row`:J$row").Style = "Normal"
